# Implement the "label and table tag" change:
#   - Add two new sheets: LabelLocators (sheetId 7) and TableLocators (sheetId 8),
#     placed after the existing LinkLocators sheet, each with the same
#     Loc1..Loc7 header row used by the other *Locators sheets.
#   - Update the selection on InputLocators (B17 -> B20) and LinkLocators
#     (E11 -> A1:G1).
#   - Make TableLocators (the last sheet) the active / selected tab.

$wb = $excel.ActiveWorkbook

$headers = "Loc1", "Loc2", "Loc3", "Loc4", "Loc5", "Loc6", "Loc7"

# --- InputLocators: selection moves from B17 to B20 ---
$wsInput = $wb.Worksheets.Item("InputLocators")
[void]$wsInput.Range("B20").Select()

# --- LinkLocators: selection moves from E11 to A1:G1 ---
$wsLink = $wb.Worksheets.Item("LinkLocators")
[void]$wsLink.Range("A1:G1").Select()

# --- New sheet: LabelLocators, right after LinkLocators ---
$wsLabel = $wb.Worksheets.Add($null, $wsLink)
$wsLabel.Name = "LabelLocators"
for ($i = 0; $i -lt $headers.Length; $i++) {
    $wsLabel.Cells.Item(1, $i + 1).Value = $headers[$i]
}
[void]$wsLabel.Range("A1:G1").Select()

# --- New sheet: TableLocators, right after LabelLocators (becomes the active tab) ---
$wsTable = $wb.Worksheets.Add($null, $wsLabel)
$wsTable.Name = "TableLocators"
for ($i = 0; $i -lt $headers.Length; $i++) {
    $wsTable.Cells.Item(1, $i + 1).Value = $headers[$i]
}
[void]$wsTable.Range("K25").Select()

# Scroll the visible tab strip so TextareaLocators (3rd sheet) is the first
# visible tab, matching the saved "firstSheet" window setting.
$win = $excel.ActiveWindow
[void]$win.ScrollWorkbookTabs($null, 3)

Write-Host "Sheets now:" ($wb.Worksheets | ForEach-Object { $_.Name })
